$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column F (dSF) - repull data, push all data, mean calculation
$updates = @{
    2  = -2
    4  = 1
    5  = -3
    6  = 7
    7  = -4
    8  = 10
    9  = 3
    11 = 2
    12 = 2
    13 = 7
    14 = -3
    15 = 1
    17 = 1
    18 = -5
    19 = 2
    20 = -6
    21 = 2
    22 = 2
    24 = -3
    25 = 1
    28 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
